# BecomePartnerPage.xlsx — add a "Locator Type" column (CSS vs Xpath) next
# to the existing Sno / Locator Name / Locator Value table, matching the
# "Execute Tests on Android mobile browser" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BecomePartnerPage")

$lastRow = 39

# --- New column D: header + per-row classification -------------------
$ws.Range("D1").Value = "Locator Type"
$ws.Range("D1").Font.Bold = $true

for ($r = 2; $r -le $lastRow; $r++) {
    $locatorValue = $ws.Cells.Item($r, 3).Value()
    if ($locatorValue -like "//*") {
        $ws.Cells.Item($r, 4).Value = "Xpath"
    } else {
        $ws.Cells.Item($r, 4).Value = "CSS"
    }
}

# --- Column widths, nudged to the dimensions left by the edit --------
$ws.Columns.Item(1).ColumnWidth = 3
$ws.Columns.Item(2).ColumnWidth = 64.5
$ws.Columns.Item(3).ColumnWidth = 69
$ws.Columns.Item(4).ColumnWidth = 31

# --- Selection follows the last touched cell --------------------------
$null = $ws.Range("D39").Select()
